$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new ticker rows after the existing data (rows 282-284)
$ws.Range("A282").Value = "IMX-USD"
$ws.Range("A283").Value = "TAO-USD"
$ws.Range("A284").Value = "MNT-USD"
